# "Haches integrees au framework" - add the Axes weapon sheet data,
# wire up the new PowerRuling entries (Lifedrain/Grounded/Axe abilities),
# and refresh the saved selections on the other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PowerRuling sheet: new ability entries used by the Axes sheet formulas
# ---------------------------------------------------------------------
$wsPR = $wb.Worksheets.Item("PowerRuling")

$wsPR.Range("B38").Value = "Axe Throw"
$wsPR.Range("C38").Value = 20

$wsPR.Range("B39").Value = "Cleave"
$wsPR.Range("C39").Value = 11

$wsPR.Range("B40").Value = "Crush"
$wsPR.Range("C40").Value = 16

$wsPR.Range("B22").Value = "Grounded"
$wsPR.Range("C22").Value = 14

$wsPR.Range("K14").Value = "Lifedrain"
$wsPR.Range("L14").Value = 10

$wsPR.Range("B41").Value = "Bloodbath"
$wsPR.Range("C41").Value = 14

$wsPR.Range("B42").Value = "Execute"
$wsPR.Range("C42").Value = 10

$wsPR.Range("B43").Value = "Endless Rampage"
$wsPR.Range("C43").Value = 16

$wsPR.Range("C44").Select()

# ---------------------------------------------------------------------
# Axes sheet: full item table
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Axes")

$ws.Columns.Item(4).ColumnWidth = 29.833333333333336
$ws.Columns.Item(5).ColumnWidth = 45.666666666666664
$ws.Columns.Item(12).ColumnWidth = 12.833333333333332

$ws.Range("F1").Value = "Health"
$ws.Range("G1").Value = "Mana"
$ws.Range("H1").Value = "Essence"
$ws.Range("I1").Value = "Strength"
$ws.Range("J1").Value = "Defense"
$ws.Range("K1").Value = "Magic"
$ws.Range("L1").Value = "Magic Defense"
$ws.Range("M1").Value = "Agility"
$ws.Range("N1").Value = "Luck"
$ws.Range("O1").Value = "Power"

# Row 2 - Iron Axe
$ws.Range("B2").Value = 1053
$ws.Range("C2").Value = "Iron Axe"
$ws.Range("D2").Value = "A simple axe, very effective to cut through many things."
$ws.Range("D2").WrapText = $true
$ws.Range("E2").Value = "Rank I"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Formula = "=F2/10+G2/5+H2+I2+J2+K2+L2+M2+N2"
$ws.Rows.Item(2).RowHeight = 36.75

# Row 3 - Steel Axe
$ws.Range("B3").Value = 1054
$ws.Range("C3").Value = "Steel Axe"
$ws.Range("D3").Value = "A simple axe, very effective to cut through many things."
$ws.Range("D3").WrapText = $true
$ws.Range("E3").Value = "Rank II, Upgrades from Iron Axe"
$ws.Range("F3").Value = 30
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 12
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Formula = "=F3/10+G3/5+H3+I3+J3+K3+L3+M3+N3"
$ws.Rows.Item(3).RowHeight = 48.75

# Row 4 - Splitting Axe
$ws.Range("B4").Value = 1055
$ws.Range("C4").Value = "Splitting Axe"
$ws.Range("D4").Value = "A simple axe, very effective to cut through many things."
$ws.Range("D4").WrapText = $true
$ws.Range("E4").Value = "Rank III, Upgrades from Steel Axe"
$ws.Range("F4").Value = 70
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 18
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Formula = "=F4/10+G4/5+H4+I4+J4+K4+L4+M4+N4"
$ws.Rows.Item(4).RowHeight = 90

# Row 5 - Sagaris
$ws.Range("B5").Value = 1056
$ws.Range("C5").Value = "Sagaris"
$ws.Range("D5").Value = "A light axe, with better precision but less strength. Inflicts Bleed on hit."
$ws.Range("D5").WrapText = $true
$ws.Range("E5").Value = "Rank IV"
$ws.Range("F5").Value = 80
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 21
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Formula = "=F5/10+G5/5+H5+I5+J5+K5+L5+M5+N5+PowerRuling!C6"
$ws.Rows.Item(5).RowHeight = 120

# Row 6 - Tomahawk
$ws.Range("B6").Value = 1057
$ws.Range("C6").Value = "Tomahawk"
$ws.Range("D6").Value = "An axe that deals heavy blows on hit. Inflicts Slow on hit."
$ws.Range("D6").WrapText = $true
$ws.Range("E6").Value = "Rank IV, Upgrades from Splitting Axe"
$ws.Range("F6").Value = 70
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 18
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Formula = "=F6/10+G6/5+H6+I6+J6+K6+L6+M6+N6+PowerRuling!C19"
$ws.Rows.Item(6).RowHeight = 75

# Row 7 - Battle Axe
$ws.Range("B7").Value = 1058
$ws.Range("C7").Value = "Battle Axe"
$ws.Range("D7").Value = "A heavy axe with a large blade. Provides additional parry."
$ws.Range("D7").WrapText = $true
$ws.Range("E7").Value = "Rank V, Upgrades from Splitting Axe"
$ws.Range("F7").Value = 80
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 26
$ws.Range("J7").Value = 4
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Formula = "=F7/10+G7/5+H7+I7+J7+K7+L7+M7+N7+1*8"
$ws.Rows.Item(7).RowHeight = 90

# Row 8 - Lumberjack's Dream
$ws.Range("B8").Value = 1059
$ws.Range("C8").Value = "Lumberjack's Dream"
$ws.Range("D8").Value = "An light axe than can cut through wood, or anything, with ease. Inflicts Bleed on hit. Increases critical strike chance."
$ws.Range("D8").WrapText = $true
$ws.Range("E8").Value = "Rank V, Upgrades from Sagaris"
$ws.Range("F8").Value = 90
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 28
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Formula = "=F8/10+G8/5+H8+I8+J8+K8+L8+M8+N8+PowerRuling!C6+1*8"
$ws.Rows.Item(8).RowHeight = 195

# Row 9 - Ono
$ws.Range("B9").Value = 1060
$ws.Range("C9").Value = "Ono"
$ws.Range("D9").Value = "A short axe with a sharp and well-defined blade. Inflicts Slow and Hemorrhage on hit."
$ws.Range("D9").WrapText = $true
$ws.Range("E9").Value = "Rank VI, Upgrades from Tomahawk"
$ws.Range("F9").Value = 110
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 34
$ws.Range("J9").Value = 6
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 6
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Formula = "=F9/10+G9/5+H9+I9+J9+K9+L9+M9+N9+PowerRuling!C19+PowerRuling!C9"
$ws.Rows.Item(9).RowHeight = 72

# Row 10 - Francisca
$ws.Range("B10").Value = 1061
$ws.Range("C10").Value = "Francisca"
$ws.Range("D10").Value = "A short axe with a long horizontal blade. Inflicts Bleed on hit. Holds the ability Axe Throw."
$ws.Range("D10").WrapText = $true
$ws.Range("E10").Value = "Rank VII, Upgrades from Sagaris"
$ws.Range("F10").Value = 150
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 43
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 7
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Formula = "=F10/10+G10/5+H10+I10+J10+K10+L10+M10+N10+PowerRuling!C6+PowerRuling!C38"
$ws.Rows.Item(10).RowHeight = 69.75

# Row 11 - The Cleaver
$ws.Range("B11").Value = 1062
$ws.Range("C11").Value = "The Cleaver"
$ws.Range("D11").Value = "A large axe with a blade capable of cutting through almost everything. Inflicts Bleed on hit. Increases critical strike chance. Holds the ability Cleave."
$ws.Range("D11").WrapText = $true
$ws.Range("E11").Value = "Rank VIII, Upgrades from Lumberjack's Dream"
$ws.Range("F11").Value = 240
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 40
$ws.Range("J11").Value = 16
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 16
$ws.Range("M11").Value = -10
$ws.Range("N11").Value = 0
$ws.Range("O11").Formula = "=F11/10+G11/5+H11+I11+J11+K11+L11+M11+N11+1*8+PowerRuling!C6+PowerRuling!C39"
$ws.Rows.Item(11).RowHeight = 120

# Row 12 - The Crusher
$ws.Range("B12").Value = 1063
$ws.Range("C12").Value = "The Crusher"
$ws.Range("D12").Value = "A heavy axe with a long haft. Its weight is often enough to heavily damage its foes, without even cutting through them. Provides additional parry. Holds the ability Crush."
$ws.Range("D12").WrapText = $true
$ws.Range("E12").Value = "Rank VIII, Upgrades from Battle Axe"
$ws.Range("F12").Value = 80
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 55
$ws.Range("J12").Value = 16
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 16
$ws.Range("M12").Value = -10
$ws.Range("N12").Value = 0
$ws.Range("O12").Formula = "=F12/10+G12/5+H12+I12+J12+K12+L12+M12+N12+1*8+PowerRuling!C40"
$ws.Rows.Item(12).RowHeight = 119.25

# Row 13 - The Sanguine
$ws.Range("B13").Value = 1064
$ws.Range("C13").Value = "The Sanguine"
$ws.Range("D13").Value = "A red-tainted axe, from the blood of its deceased foes. Holds the ability Bloodbath. Heals for a portion of damage dealt."
$ws.Range("D13").WrapText = $true
$ws.Range("E13").Value = "Rank VIII"
$ws.Range("F13").Value = 200
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 48
$ws.Range("J13").Value = 14
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 14
$ws.Range("M13").Value = -10
$ws.Range("N13").Value = 0
$ws.Range("O13").Formula = "=F13/10+G13/5+H13+I13+J13+K13+L13+M13+N13+PowerRuling!C41+PowerRuling!L14"
$ws.Rows.Item(13).RowHeight = 75.75

# Row 14 - Parashu
$ws.Range("B14").Value = 1065
$ws.Range("C14").Value = "Parashu"
$ws.Range("D14").Value = "A short axe made of rare ore. It makes it light-weighted. Inflicts Bleed on hit. Holds the ability Axe Throw. Increases Agility and Precision."
$ws.Range("D14").WrapText = $true
$ws.Range("E14").Value = "Rank IX, Upgrades from Francisca"
$ws.Range("F14").Value = 200
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 49
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 20
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Formula = "=F14/10+G14/5+H14+I14+J14+K14+L14+M14+N14+PowerRuling!C6+PowerRuling!C38+2*8"
$ws.Rows.Item(14).RowHeight = 96.75

# Row 15 - Yue
$ws.Range("B15").Value = 1066
$ws.Range("C15").Value = "Yue"
$ws.Range("D15").Value = "A long haft with a small but very sharp blade. Inflicts Slow and Hemorrhage on hit. Holds the ability Execute."
$ws.Range("D15").WrapText = $true
$ws.Range("E15").Value = "Rank IX, Upgrades from Ono"
$ws.Range("F15").Value = 200
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 55
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 15
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Formula = "=F15/10+G15/5+H15+I15+J15+K15+L15+M15+N15+PowerRuling!C9+PowerRuling!C19+PowerRuling!C42"
$ws.Rows.Item(15).RowHeight = 91.5

# Row 16 - Golden Axe
$ws.Range("B16").Value = 1067
$ws.Range("C16").Value = "Golden Axe"
$ws.Range("D16").Value = "A heavy axe made of precious ore that prevents its blade from deteriorating. Inflicts Bleed on hit. Holds the abilities Axe Throw and Endless Rampage. Increases Agility and Precision."
$ws.Range("D16").WrapText = $true
$ws.Range("E16").Value = "Rank X, Upgrades from Parashu"
$ws.Range("F16").Value = 200
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 49
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 24
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Formula = "=F16/10+G16/5+H16+I16+J16+K16+L16+M16+N16+PowerRuling!C6+2*8+PowerRuling!C38+PowerRuling!C43"
$ws.Rows.Item(16).RowHeight = 114.75

# Row 17 - Warchief's Axe
$ws.Range("B17").Value = 1068
$ws.Range("C17").Value = "Warchief's Axe"
$ws.Range("D17").Value = "A large axe that sweeps everything in its path. Provides additional parry and critical strike chance. Inflicts Bleed on hit. Holds the abilities Cleave and Crush."
$ws.Range("D17").WrapText = $true
$ws.Range("E17").Value = "Rank XI, Upgrades from The Cleaver + The Crusher"
$ws.Range("F17").Value = 300
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 65
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = -10
$ws.Range("N17").Value = 0
$ws.Range("O17").Formula = "=F17/10+G17/5+H17+I17+J17+K17+L17+M17+N17+2*8+PowerRuling!C6+PowerRuling!C39+PowerRuling!C40"
$ws.Rows.Item(17).RowHeight = 110.25

# Row 18 - Warlord's Steel
$ws.Range("B18").Value = 1069
$ws.Range("C18").Value = "Warlord's Steel"
$ws.Range("D18").Value = "A large, red-tainted blade, testimony of its wielder's glorious past. Provides additional Parry. Heals for a portion of damage dealt. Holds the abilities Bloodbath and Crush."
$ws.Range("D18").WrapText = $true
$ws.Range("E18").Value = "Rank XI, Upgrades from the Crusher + The Sanguine"
$ws.Range("F18").Value = 200
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 80
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 18
$ws.Range("M18").Value = -10
$ws.Range("N18").Value = 0
$ws.Range("O18").Formula = "=F18/10+G18/5+H18+I18+J18+K18+L18+M18+N18+1*8+PowerRuling!L14+PowerRuling!C40+PowerRuling!C41"
$ws.Rows.Item(18).RowHeight = 100.5

# Row 19 - Blood Angel
$ws.Range("B19").Value = 1070
$ws.Range("C19").Value = "Blood Angel"
$ws.Range("D19").Value = "A gigantic blade, that brings his foes closer to another realm."
$ws.Range("D19").WrapText = $true
$ws.Range("E19").Value = "Rank XII, Warlok's ultimate weapon"
$ws.Range("O19").Formula = "=F19/10+G19/5+H19+I19+J19+K19+L19+M19+N19"
$ws.Rows.Item(19).RowHeight = 54.75

$ws.Activate()
$ws.Range("I18").Select()

# ---------------------------------------------------------------------
# Refresh saved selections on the other sheets (scroll/selection state)
# ---------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("ListOfItems")
$wsItems.Activate()
$wsItems.Range("B67:F84").Select()

$wsSwords = $wb.Worksheets.Item("Swords")
$wsSwords.Activate()
$wsSwords.Range("M34").Select()

$wsDaggers = $wb.Worksheets.Item("Daggers")
$wsDaggers.Activate()
$wsDaggers.Range("F1:O1").Select()

$ws.Activate()
$ws.Range("I18").Select()
